$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row via Ctrl+Down-style navigation from the header.
$lastRow = $ws.Cells.Item(1, 1).End(4).Row

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45177) {
        $cell.Value2 = 45178
    }
}
